$d = $word.ActiveDocument

# Locate the "Professor: Chandra" text and collapse the found range to its
# end, so we can insert " Bobba" immediately after "Chandra" (before the
# following tab character).
$r = $d.Content
$found = $r.Find.Execute("Professor: Chandra", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $r.Collapse(0)
    $r.InsertAfter(" Bobba")
}
